# Fix minor errors referencing "Person" instead of "Student" in the UI
# component class diagram:
#   - "PersonListPanel" rectangle label -> "StudentListPanel"
#     (also nudge its font down from 10.5pt to 10pt so the longer
#     label still fits inside the same box)
#   - "PersonCard" rectangle label -> "StudentCard"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$listPanel = $null
$card = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "PersonListPanel") {
            $listPanel = $shp
        } elseif ($txt -eq "PersonCard") {
            $card = $shp
        }
    }
}

# Fall back to the known positions on this slide if the text-based
# lookup above didn't find them for some reason.
if ($listPanel -eq $null) {
    $listPanel = $s.Shapes.Item(11)
}
if ($card -eq $null) {
    $card = $s.Shapes.Item(12)
}

$listPanelRange = $listPanel.TextFrame.TextRange
$listPanelRange.Text = "StudentListPanel"
$listPanelRange.Font.Size = 10

$card.TextFrame.TextRange.Text = "StudentCard"
